# Fruta / hortaliza, semanal
# Insert a new weekly record as row 69 in the "Pomelo" price sheet.
# This pushes the existing rows 69-152 down to 70-153 (dimension grows
# from A1:T152 to A1:T153) and populates the newly-opened row 69 with
# the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 69, shifting rows 69:152 down
# to 70:153 (and extending the sheet's used range to T153).
$ws.Rows.Item(69).Insert()

# Populate the newly inserted row 69 with this week's record.
$ws.Cells.Item(69, 1).Value  = 4
$ws.Cells.Item(69, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(69, 3).Value  = "Los Lagos"
$ws.Cells.Item(69, 4).Value  = 44483
$ws.Cells.Item(69, 5).Value  = 10
$ws.Cells.Item(69, 6).Value  = "Fruta"
$ws.Cells.Item(69, 7).Value  = 100102
$ws.Cells.Item(69, 8).Value  = "Cítricos"
$ws.Cells.Item(69, 9).Value  = 100102006
$ws.Cells.Item(69, 10).Value = "Pomelo"
$ws.Cells.Item(69, 11).Value = "Start Ruby"
$ws.Cells.Item(69, 12).Value = "Primera"
$ws.Cells.Item(69, 13).Value = 80
$ws.Cells.Item(69, 14).Value = 11000
$ws.Cells.Item(69, 15).Value = 12000
$ws.Cells.Item(69, 16).Value = 11500
$ws.Cells.Item(69, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(69, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(69, 19).Value = 821
$ws.Cells.Item(69, 20).Value = 14
